# ICPP: Make pack image numbers easier to read.
# Bump the font size (20pt) and bold the single-digit number labels on the
# four "pack" rectangles that were still using the default run formatting.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# These are the shape Ids (p:cNvPr id=) of the four number labels touched by
# the commit - identified by their unique ids rather than positional index
# so the script is resilient to shape reordering.
$targetIds = @(29, 49, 57, 58)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($targetIds -contains $shp.Id) {
        $tr = $shp.TextFrame.TextRange
        $tr.Font.Size = 20
        $tr.Font.Bold = $true
    }
}
